# Auto-applied edits per xml_diff for cryptos.xlsx crypto-list update
# (GitHub Actions hourly refresh of coinranking.com snapshot)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value updates (price / volume% columns) ---
# Cells whose new text happens to parse as a plain number (e.g. "586.82")
# are forced back to Text so Excel doesn't silently convert them to a
# numeric cell the way plain .Value assignment would; the style index is
# then restored to Normal so no stray formatting is introduced.
$ws.Range("D2").Value = '67.761.44'
$ws.Range("E2").Value = '  +3.00%  '
$ws.Range("D3").Value = '3.327.05'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +3.44%  '
$ws.Range("D9").Value = '3.320.80'
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("E10").Value = '  +3.24%  '
$ws.Range("E11").Value = '  +1.43%  '
$ws.Range("E12").Value = '  +2.03%  '
$ws.Range("E13").Value = '  +5.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '638.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +10.93%  '
$ws.Range("D15").Value = '3.858.74'
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("D17").Value = '67.834.95'
$ws.Range("E17").Value = '  +3.35%  '
$ws.Range("E18").Value = '  +1.79%  '
$ws.Range("D19").Value = '3.327.10'
$ws.Range("E19").Value = '  +0.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.25%  '
$ws.Range("E21").Value = '  +1.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.898'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.63'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("E26").Value = '  +2.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.25%  '
$ws.Range("E28").Value = '  +3.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +8.03%  '
$ws.Range("E30").Value = '  +3.12%  '
$ws.Range("E31").Value = '  +1.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '591.81'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.97%  '
$ws.Range("D33").Value = '3.936.73'
$ws.Range("E33").Value = '  +5.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '10.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.53'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.66%  '
$ws.Range("E36").Value = '  +2.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.70'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.43%  '
$ws.Range("E39").Value = '  +4.89%  '
$ws.Range("E40").Value = '  +1.69%  '
$ws.Range("E41").Value = '  +4.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '32.65'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.09%  '
$ws.Range("E45").Value = '  +2.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0414'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.69%  '
$ws.Range("E47").Value = '  +2.75%  '
$ws.Range("E48").Value = '  +0.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.54'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.59%  '
$ws.Range("E50").Value = '  +9.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '130.59'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.92%  '

# --- Row 43/44 swap: PEPE and ApeXProtocol switch rank order ---
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.38'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.68%  '
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '0.0₃0684'
# E44 stays "  +1.74%  " (unchanged per diff)
